# Rename existing sheet and rebuild it as the "rule_desc" description table,
# then add a new "id_1" sheet holding the assessment weightage rules (as an
# Excel Table), matching the "New files with flask application" commit.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "rule_desc"

# ---- Sheet 1: rule_desc ----
$ws1.Range("A1").Value = "id"
$ws1.Range("B1").Value = "desc"
$ws1.Range("C1").Value = "rule_func()"

$ws1.Range("A2").Value = 1
$ws1.Range("B2").Value = "Displays if any of the assessment is breaching its expected weightage in the syllabus"
$ws1.Range("C2").Value = "rule_1(page1)"

$ws1.Range("A3").Value = 2
$ws1.Range("B3").Value = "Displays verb distance of EKS and CLO based on cognitive domain word list"
$ws1.Range("C3").Value = "rule_2_cd(page1)"

$ws1.Range("A4").Value = 3
$ws1.Range("B4").Value = "Displays verb distance of EKS and CLO based on psychomotor domain word list"
$ws1.Range("C4").Value = "rule_2_pd(page1)"

$ws1.Range("A5").Value = 4
$ws1.Range("B5").Value = "Displays verb distance of EKS and CLO based on affective domain word list"
$ws1.Range("C5").Value = "rule_2_ad(page1)"

$ws1.Columns.Item(2).ColumnWidth = 54.88671875
$ws1.Columns.Item(3).ColumnWidth = 15.21875

$ws1.Range("B1").Select() | Out-Null

# ---- Sheet 2: id_1 (new sheet, inserted right after rule_desc) ----
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "id_1"

$ws2.Range("A1").Value = "assessment"
$ws2.Range("B1").Value = "lower"
$ws2.Range("C1").Value = "higher"
$ws2.Range("D1").Value = "rules"

$ws2.Range("A2").Value = "Quiz"
$ws2.Range("B2").Value = 1
$ws2.Range("C2").Value = 4
$ws2.Range("D2").Value = "Quiz weightage should be between 1% to 4%"

$ws2.Range("A3").Value = "Test"
$ws2.Range("B3").Value = 5
$ws2.Range("C3").Value = 24
$ws2.Range("D3").Value = "Test weightage should be between 5% to 24%"

$ws2.Range("A4").Value = "Exam"
$ws2.Range("B4").Value = 25
$ws2.Range("C4").Value = 40
$ws2.Range("D4").Value = "Exam weightage should be between 25% to 40%"

$ws2.Columns.Item(1).ColumnWidth = 12.44140625
$ws2.Columns.Item(4).ColumnWidth = 40.77734375

# Turn the range into a proper Excel Table ("Table3"), no AutoFilter arrows shown.
$tbl = $ws2.ListObjects.Add(1, $ws2.Range("A1:D4"), $null, 1)
$tbl.Name = "Table3"
$tbl.ShowAutoFilter = $false

# id_1 is the tab that ends up active/selected.
$ws2.Activate() | Out-Null
